$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "FEE"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Name"

$ws.Range("A2").Value = 111
$ws.Range("B2").Value = 500
$ws.Range("D2").Value = "Tajinder "

$ws.Range("A3").Value = 115
$ws.Range("B3").Value = 500
$ws.Range("D3").Value = "Prabh"

$ws.Range("A4").Value = 236
$ws.Range("B4").Value = 200
$ws.Range("C4").Value = "Concession Applied"
$ws.Range("D4").Value = "Surajpal"

$ws.Range("A5").Value = 536
$ws.Range("B5").Value = 500
$ws.Range("D5").Value = "Manjot"

$ws.Range("D12").Select() | Out-Null
